# Reorders the comma-separated "Recorded By" names in column G so that
# "System" (exact case) sorts first, followed by other non-email names,
# followed by email addresses (case-insensitive alphabetical within each
# group). Cells with a single value (no comma) are left untouched.
#
# Note: Sort-Object with a scriptblock -Property, or with multiple
# -Property keys, is not reliable in this runtime, so we build helper
# objects with precomputed sort keys and chain single-key stable sorts
# (secondary key first, then primary key) to get a correct multi-key sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $raw = $cell.Value2

    if ($null -eq $raw) { continue }
    if (-not ($raw -is [string])) { continue }
    if ($raw.IndexOf(",") -lt 0) { continue }

    $parts = $raw -split "," | ForEach-Object { $_.Trim() }

    $objs = @()
    foreach ($p in $parts) {
        $rank = 2
        if ($p.Equals("System")) { $rank = 0 }
        elseif ($p.IndexOf("@") -lt 0) { $rank = 1 }
        $objs += [PSCustomObject]@{ Name = $p; Rank = $rank; Key = $p.ToLower() }
    }

    $byKey = $objs | Sort-Object -Property Key
    $byRank = $byKey | Sort-Object -Property Rank

    $names = $byRank | ForEach-Object { $_.Name }
    $newVal = [string]::Join(", ", $names)

    # NOTE: -eq/-ne/-ceq/-cne are all case-insensitive in this runtime,
    # so use .Equals() (case-sensitive) to decide whether an update is needed.
    if (-not $newVal.Equals($raw)) {
        $cell.Value = $newVal
    }
}
